$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "57.799.19"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "2.420.22"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "510.37"
$ws.Range("E5").Value = "  -0.57%  "
Set-TextValue "D6" "133.32"
$ws.Range("E6").Value = "  +5.18%  "
Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue "D8" "0.557"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").Value = "2.445.92"
$ws.Range("E9").Value = "  +3.53%  "
Set-TextValue "D10" "0.0975"
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +2.68%  "
Set-TextValue "D13" "4.62"
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").Value = "2.869.71"
$ws.Range("E14").Value = "  +3.80%  "
$ws.Range("D15").Value = "57.448.53"
$ws.Range("E15").Value = "  +1.80%  "
Set-TextValue "D16" "21.99"
$ws.Range("E16").Value = "  +3.84%  "
Set-TextValue "D17" "0.0000134"
$ws.Range("E17").Value = "  +4.03%  "
$ws.Range("D18").Value = "2.492.36"
$ws.Range("E18").Value = "  +5.66%  "
Set-TextValue "D19" "10.32"
$ws.Range("E19").Value = "  +1.36%  "
Set-TextValue "D20" "314.75"
$ws.Range("E20").Value = "  +2.34%  "
Set-TextValue "D21" "4.08"
$ws.Range("E21").Value = "  +2.40%  "
Set-TextValue "D22" "6.39"
$ws.Range("E22").Value = "  +6.57%  "
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "5.79"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D24" "0.998"
$ws.Range("E24").Value = "  +0.11%  "
Set-TextValue "D25" "66.05"
$ws.Range("E25").Value = "  +2.58%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  +0.12%  "
Set-TextValue "D27" "0.155"
$ws.Range("E27").Value = "  +1.34%  "
Set-TextValue "D28" "0.383"
$ws.Range("E28").Value = "  -1.83%  "
Set-TextValue "D29" "7.58"
$ws.Range("E29").Value = "  +6.36%  "
Set-TextValue "D30" "170.39"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "0.0₃0733"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  +2.72%  "
Set-TextValue "D33" "6.04"
$ws.Range("E33").Value = "  -0.51%  "
Set-TextValue "D34" "1.13"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +0.12%  "
Set-TextValue "D36" "0.996"
$ws.Range("E36").Value = "  +0.10%  "
Set-TextValue "D37" "18.06"
$ws.Range("E37").Value = "  +3.40%  "
Set-TextValue "D38" "1.21"
$ws.Range("E38").Value = "  +3.82%  "
Set-TextValue "D39" "3.87"
$ws.Range("E39").Value = "  +5.39%  "
Set-TextValue "D40" "36.76"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D41" "1.46"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D42" "0.805"
$ws.Range("E42").Value = "  +1.67%  "
Set-TextValue "D43" "132.99"
$ws.Range("E43").Value = "  +9.24%  "
$ws.Range("E44").Value = "  +3.36%  "
Set-TextValue "D45" "4.99"
$ws.Range("E45").Value = "  +3.32%  "
Set-TextValue "D46" "254.99"
$ws.Range("E46").Value = "  +2.27%  "
Set-TextValue "D47" "0.571"
$ws.Range("E47").Value = "  +1.29%  "
Set-TextValue "D48" "0.0914"
$ws.Range("E48").Value = "  +1.38%  "
Set-TextValue "D49" "0.0492"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "17.33"
$ws.Range("E50").Value = "  +5.67%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D51" "0.0212"
$ws.Range("E51").Value = "  +3.44%  "
